$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sedan_HambaLG")

# --- Updated hardpoint values (row 5 / row 6 / row 9 / row 10) ---
# Row 5 (Front / sTop)
$ws.Range("G5").Value = 0.62
$ws.Range("H5").Value = 0.65
$ws.Range("G5:H5").NumberFormat = "0.00"

# Row 6 (Front / sBottom)
$ws.Range("G6").Value = 0.85
$ws.Range("H6").Value = 0.19
$ws.Range("G6:H6").NumberFormat = "0.00"

# Row 9 (Rear / sTop) - F9 keeps its value, but number format changes too
$ws.Range("G9").Value = 0.62
$ws.Range("H9").Value = 0.65
$ws.Range("F9:H9").NumberFormat = "0.00"

# Row 10 (Rear / sBottom) - F10 keeps its value, but number format changes too
$ws.Range("G10").Value = 0.85
$ws.Range("H10").Value = 0.19
$ws.Range("F10:H10").NumberFormat = "0.00"

# --- Cosmetic: sheet tab color (theme 9 -> theme 8) ---
$ws.Tab.Color = 11957550
